# Apply updated dSF (column F) values per row, per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 1
    14 = 2
    15 = 0
    16 = 2
    22 = 1
    27 = 1
    28 = 3
    29 = -3
    45 = 0
    46 = 0
    53 = 1
    57 = 2
    58 = 1
    59 = 0
    60 = -1
    62 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
